$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-0.8474283477698586, 1.1311697144454378, 1.1302973865769903, 3.0681335472858438)"
$ws.Range("C2").Value = "JSU(-1.4149470145952099, 1.3154210372555823, 2.122257604719718, 4.998791519567386)"
$ws.Range("D2").Value = "NIG(0.6723856071459912, 0.4974495385655613, 0.8987594262835561, 2.8429787326487426)"
$ws.Range("E2").Value = "JSU(-1.5732717628558524, 1.1913249791836773, 1.4266563860350547, 4.519897210492166)"
